$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-21 (columns A,B,C,D,E)
$data = @(
    @(20,  2, 40, "3",        "25"),
    @(33,  3, 50, "5,8",      "19,19"),
    @(37,  3, 50, "6,15,21",  "18,24,2"),
    @(47,  3, 50, "14,21",    "28,22"),
    @(57,  1, 30, "2",        "26"),
    @(59,  3, 50, "4,10",     "27,23"),
    @(74,  3, 50, "1,7",      "28,22"),
    @(153, 3, 50, "13,20,21", "36,12,2"),
    @(240, 3, 50, "10,12",    "4,46"),
    @(248, 2, 40, "12,18",    "1,19"),
    @(257, 1, 30, "7,9",      "3,27"),
    @(314, 1, 30, "9",        "3"),
    @(378, 3, 50, "11,16",    "24,26"),
    @(442, 3, 50, "17",       "35"),
    @(470, 3, 50, "27",       "47"),
    @(478, 3, 50, "19,25,26", "35,2,7"),
    @(482, 1, 30, "16",       "1"),
    @(509, 3, 50, "23,24",    "18,11"),
    @(549, 2, 40, "22,24",    "18,22"),
    @(580, 2, 40, "25",       "40")
)

# Columns D and E hold text labels (lists of community/people ids), even
# when the text happens to look like a plain integer. Force text format
# first so Excel doesn't silently reinterpret them as numbers.
$ws.Range("D2:E21").NumberFormat = "@"

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $rowIndex++
}

# Remove the now-obsolete last row (row 22), shrinking the sheet to A1:E21
$ws.Range("A22:E22").Delete()
